$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CasesTab" query (cell B2) had a stray, broken `Cohort` column in its
# RETURN clause (the commit "Fixed variables and query errors ... from TC01
# to TC30"). Drop that trailing RETURN item while leaving the rest of the
# Cypher text (including the now-unused `co` match/WITH variable) untouched.
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN['German Shepherd Dog'] 
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $casesQuery

# Move the saved selection/active cell from B4 to B2 (matches the workbook's
# last-saved cursor position recorded in the sheet view).
$ws.Range("B2").Select()
